# "Fix Enum val is empty" - the "Type" / "int Enum[EItemType]" / "类型" / "武器"
# column (column D: Type/EItemType) is empty of real data (duplicate of the
# SubType enum) and gets removed entirely from the ItemData sheet. Deleting
# the whole column shifts the old Quality/SubType columns (E,F) left into
# D,E, and Excel automatically re-indexes sheetData, column widths,
# dataValidations, dimension and the shared-string table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemData")

$ws.Columns("D:D").Delete()

# The conditional-formatting rules that colour the quality/sub-type columns
# ("containsText" on 紫色/灰色/绿色) referenced the old column E; after the
# column shift they need to point at the new column D instead.
$fc = $ws.Cells.FormatConditions
$newRange = $ws.Range("D1:E1048576")
$colors = @("紫色", "灰色", "绿色")
for ($i = 2; $i -le $fc.Count; $i++) {
    $cond = $fc.Item($i)
    $cond.ModifyAppliesToRange($newRange)
    $color = $colors[$i - 2]
    $cond.Formula1 = '=NOT(ISERROR(SEARCH("' + $color + '",D1)))'
}

# Restore the user's on-screen selection after the edit.
[void]$ws.Range("I11").Select()
